# Applies the "Add row with correct predicates" edit:
#  - header sheet: fix the "nesp: TBA" placeholder to a real curie URL and
#    add a new "sssom:" curie-map line
#  - SSSOM sheet: add a new RDF/SSSOM predicate header row above the existing
#    header row (pushing the existing header + data rows down by one), clean
#    up the nesp curie labels in column A, and add a helper "rdfs:label"
#    formula column (O) to every data row

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "header" sheet
# ---------------------------------------------------------------------------
$header = $wb.Worksheets.Item("header")

# Row 3 held a placeholder; replace it with the real nesp curie definition.
$header.Range("A3").Value = "   nesp: https://w3id.org/env/neap/nesp/"

# Insert a new row for the sssom: curie, between the orcid (row7) and
# status (row8) lines, and give it the "Normal 2" cell-style font that was
# pasted in from elsewhere.
$header.Rows.Item(8).Insert()
$header.Range("A8").Value = "   sssom: https://w3id.org/sssom/"
$header.Range("A8").Font.Name = "Calibri"
$header.Range("A8").Font.Size = 11
$header.Rows.Item(8).Select()

# ---------------------------------------------------------------------------
# 2. "SSSOM" sheet
# ---------------------------------------------------------------------------
$sssom = $wb.Worksheets.Item("SSSOM")

# Normalise the nesp curie labels used in column A (drop parentheses/& so
# they are valid curies).
$sssom.Range("A2").Value = "nesp:On-shelf-neritic-epinesp"
$sssom.Range("A3").Value = "nesp:Off-shelf-oceanic-epinesp"
$sssom.Range("A5").Value = "nesp:Bathynesp-n-Abyssonesp"
$sssom.Range("A6").Value = "nesp:Bathynesp-n-Abyssonesp"

# Insert a new top header row describing the RDF/SSSOM predicates used by
# each column of the existing header row. This pushes the existing header
# and all five data rows down by one (and the trailing blank row becomes
# row 8).
$sssom.Rows.Item(1).Insert()

$sssom.Range("A1").Value = "rdf:subject"
$sssom.Range("B1").Value = "sssom:subject_label"
$sssom.Range("C1").Value = "rdf:predicate"
$sssom.Range("D1").Value = "rdf:object"
$sssom.Range("E1").Value = "sssom:object_label"
$sssom.Range("F1").Value = "sssom:mapping_justification"
$sssom.Range("G1").Value = "dcterms:creator"
$sssom.Range("H1").Value = "sssom:creator_label"
$sssom.Range("I1").Value = "dcterms:created"
$sssom.Range("J1").Value = "sssom:confidence"
$sssom.Range("K1").Value = "crosswalk:status"
$sssom.Range("L1").Value = "sssom:reviewer_id"
$sssom.Range("M1").Value = "sssom:reviewer_id"
$sssom.Range("N1").Value = "rdfs:comment"
$sssom.Range("O1").Value = "rdfs:label"

$sssom.Range("A1:O1").Font.Bold = $true
$sssom.Range("A1:O1").Font.Name = "Calibri"
$sssom.Range("A1:O1").Font.Size = 9
$sssom.Range("A1:O1").HorizontalAlignment = -4131 # xlLeft
$sssom.Range("A1:O1").VerticalAlignment = -4108   # xlCenter
$sssom.Rows.Item(1).RowHeight = 12.75

# Add the "author_label" helper formula to every data row (3-7), mirroring
# the style already used on column G of that row.
for ($r = 3; $r -le 7; $r++) {
    $cell = $sssom.Range("O$r")
    $cell.Formula = "=CONCAT(B$r, `" - mapping to IUCN GET`")"
    $cell.Style = $sssom.Range("G$r").Style
}

# Leave SSSOM as the active sheet/selection, scrolled over so column O (the
# new formula) is visible.
$sssom.Activate()
$sssom.Application.ActiveWindow.ScrollColumn = 8
$sssom.Range("O3").Select()

$wb.Save()
